# [Draw schematic for the display unit #8]
# committing pin assignment for N5150M8CD
#
# Adds a new "Alternative" pin-name column (G) to the pin assignment
# table, giving the N5150M8CD module pin for every signal that already
# has a BMD-200 pin in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "N5150M8CD"
$ws.Range("G5").Value = "B4"
$ws.Range("G6").Value = "B3"
$ws.Range("G7").Value = "F2"
$ws.Range("G8").Value = "E3"
$ws.Range("G9").Value = "F4"
$ws.Range("G10").Value = "F3"
$ws.Range("G11").Value = "E4"
$ws.Range("G12").Value = "F5"

$ws.Columns.Item(7).AutoFit() | Out-Null

$ws.Range("G13").Select() | Out-Null
